$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (same cell, new time) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 02:22"

# --- Estados Unidos (row 4): refreshed totals ---
$ws.Range("B4").Value = 366906
$ws.Range("C4").Value = 30233
$ws.Range("D4").Value = 19656
$ws.Range("E4").Value = 336382
$ws.Range("G4").Value = 1252
$ws.Range("H4").Value = 10868

# --- Canada (row 16): refreshed totals ---
$ws.Range("B16").Value = 16667
$ws.Range("C16").Value = 1155
$ws.Range("D16").Value = 3616
$ws.Range("E16").Value = 12728

# --- Australia (row 24): refreshed totals ---
$ws.Range("E24").Value = 3418
$ws.Range("G24").Value = 8
$ws.Range("H24").Value = 45

# --- Argentina overtakes Colombia in the ranking, inserting itself ---
# --- right after Sudafrica (row 51) and pushing Colombia/Islandia down ---
# Row 52: now Argentina, with freshly updated totals
$ws.Range("A52").Value = "Argentina"
$ws.Range("B52").Value = 1628
$ws.Range("C52").Value = 74
$ws.Range("D52").Value = 325
$ws.Range("E52").Value = 1250
$ws.Range("F52").Value = 94
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 53

# Row 53: now Colombia (previous row-52 figures)
$ws.Range("A53").Value = "Colombia"
$ws.Range("B53").Value = 1579
$ws.Range("C53").Value = 94
$ws.Range("D53").Value = 88
$ws.Range("E53").Value = 1445
$ws.Range("F53").Value = 50
$ws.Range("G53").Value = 11
$ws.Range("H53").Value = 46

# Row 54: now Islandia (previous row-53 figures)
$ws.Range("A54").Value = "Islandia"
$ws.Range("B54").Value = 1562
$ws.Range("C54").Value = 76
$ws.Range("D54").Value = 460
$ws.Range("E54").Value = 1096
$ws.Range("F54").Value = 11
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 6

# Row 55 (Argelia) is unchanged.

# --- Jamaica (row 138): refreshed totals ---
$ws.Range("B138").Value = 59
$ws.Range("C138").Value = 1
$ws.Range("E138").Value = 48
